$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: the phone number in A3 was previously stored as text; convert it
# to a genuine numeric value (76442781), matching the rest of the column.
$ws.Range("A3").Value = 76442781

# Row 4: new payment record for 76442781 (Check), 2025-08-20T08:55:01
# Column A must stay text (it looked like this before the A3 conversion),
# even though its content is all digits, so force a text format, assign
# it, then clear the explicit formatting back off the cell.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "76442781"
$ws.Range("A4").ClearFormats()

$ws.Range("B4").Value = 408
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 40
$ws.Range("E4").Value = 500
$ws.Range("F4").Value = 368
$ws.Range("G4").Value = "Check"
$ws.Range("H4").Value = "2025-08-20T08:55:01"
